$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF) - match the existing header style
# (bold font, thin border all around, centered horizontal / top vertical alignment)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# New data columns I (I0) and J (IF) for rows 2-22
$data = @(
    @(4,5),
    @(7,7),
    @(6,7),
    @(6,6),
    @(7,8),
    @(9,9),
    @(7,8),
    @(5,6),
    @(7,8),
    @(6,6),
    @(8,8),
    @(6,7),
    @(5,6),
    @(7,7),
    @(7,7),
    @(9,9),
    @(7,7),
    @(3,3),
    @(3,3),
    @(6,6),
    @(8,8)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $data[$r][0]
    $ws.Cells.Item($row, 10).Value = $data[$r][1]
}
